# The document contains five "<id>...</id>" tag paragraphs that were each
# split across three separate runs:
#   run 1: "<id>"            (Courier New, color 7f6000, sz 18)
#   run 2: "<the actual id>" (plain/default formatting)
#   run 3: "</id>"           (Courier New, color 7f6000, sz 18)
#
# The edit collapses each triple into a single run containing the full
# "<id>...</id>" text, using the Courier New "tag" formatting throughout.
# Doing the replacement with Find/Execute merges the matched runs into one,
# picking up the formatting of the first run in the match - exactly what we
# want here.

$d = $word.ActiveDocument

$ids = @("p030r_2", "p031r_1", "p031r_2", "p031r_3", "p031r_4")

foreach ($id in $ids) {
    $tag = "<id>" + $id + "</id>"
    $found = $d.Content.Find.Execute($tag, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $tag, 2)
    Write-Output "merged $tag -> $found"
}
